$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from serial date 45185 to 45204 for rows 2-101
for ($row = 2; $row -le 101; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
